$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 2007 data row (row 2) entirely; rows below shift up.
$ws.Rows.Item(2).Delete()
